$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "Done By:" header to "Done By"
$ws.Range("N2").Value = "Done By"

# 2. Move the task cards (currently in the merged "Needs Reviewing" column F:G)
#    into the merged "Done" column H:I, for rows 3-7 (id 1..5 tasks all reviewed/accepted).
$ws.Range("F3:G3").Cut($ws.Range("H3:I3"))
$ws.Range("F4:G4").Cut($ws.Range("H4:I4"))
$ws.Range("F5:G5").Cut($ws.Range("H5:I5"))
$ws.Range("F6:G6").Cut($ws.Range("H6:I6"))
$ws.Range("F7:G7").Cut($ws.Range("H7:I7"))

# 3. Resize columns: "Needs Reviewing" (G) shrinks, "Done By" area (J:K) widens.
$ws.Columns.Item(7).ColumnWidth = 16.29
$ws.Range("J:K").ColumnWidth = 16.86
